$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New timeline entries (Chapter 4 and 5)
$rows = @(
    @{ A = -1800; B = $null; C = "Early trace of gold in Karnataka." },
    @{ A = -1000; B = $null; C = "Earliest date of the Khetri copper belts." },
    @{ A = 100;   B = $null; C = "Discovery of the direction of the monsoon." },
    @{ A = -500;  B = $null; C = "Wide use of iron tools in the Gangetic plains and spurt in settlements." },
    @{ A = -300;  B = $null; C = "Famine and Jain migration from Magadha to south India." },
    @{ A = 1600;  B = 1700;  C = "Forests in the Ganga-Yamuna doab despite deforestation." },
    @{ A = 1869;  B = $null; C = "The term Ecology coined." }
)

$startRow = 35
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $entry = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $entry.A
    if ($null -ne $entry.B) {
        $ws.Cells.Item($r, 2).Value = $entry.B
    }
    $ws.Cells.Item($r, 3).Value = $entry.C
}

# Scroll/select like the final saved state
$ws.Application.Goto($ws.Range("A25"), $true)
$ws.Range("C42").Select()
